# Update "last modified" timestamps (column E) for a set of tables
# to reflect the latest auto-commit regeneration pass.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")

$updates = @{
    8 = "2022年08月24日 10:44:22"
    32 = "2022年08月26日 11:11:27"
    57 = "2022年08月26日 10:22:55"
    75 = "2022年08月26日 17:12:56"
    82 = "2022年08月19日 14:25:34"
    104 = "2022年08月19日 10:44:25"
    137 = "2022年08月26日 11:42:12"
    223 = "2022年08月19日 15:29:23"
    225 = "2022年08月19日 15:33:04"
    227 = "2022年08月19日 16:17:21"
    229 = "2022年08月19日 15:31:04"
    231 = "2022年08月19日 15:31:26"
    233 = "2022年08月19日 15:31:50"
    235 = "2022年08月19日 15:32:11"
    237 = "2022年08月19日 15:32:26"
    239 = "2022年08月19日 15:32:44"
    241 = "2022年08月19日 15:33:32"
    243 = "2022年08月19日 15:34:02"
    245 = "2022年08月19日 15:34:12"
    247 = "2022年08月19日 15:34:29"
    249 = "2022年08月19日 15:34:48"
    251 = "2022年08月19日 15:35:11"
    253 = "2022年08月19日 15:35:31"
    255 = "2022年08月19日 15:35:46"
    257 = "2022年08月19日 15:35:59"
    259 = "2022年08月19日 15:36:14"
    261 = "2022年08月19日 15:36:26"
    263 = "2022年08月19日 15:36:38"
    265 = "2022年08月19日 15:36:49"
    267 = "2022年08月19日 15:36:57"
    269 = "2022年08月19日 15:37:05"
    271 = "2022年08月19日 15:37:14"
    273 = "2022年08月19日 15:37:23"
    275 = "2022年08月19日 15:37:31"
    277 = "2022年08月19日 15:37:39"
    279 = "2022年08月19日 15:37:49"
    281 = "2022年08月19日 15:37:58"
    283 = "2022年08月19日 15:38:17"
    285 = "2022年08月19日 15:38:27"
    287 = "2022年08月19日 15:38:35"
    289 = "2022年08月19日 15:38:43"
    291 = "2022年08月19日 15:38:51"
    293 = "2022年08月19日 15:38:58"
    296 = "2022年08月25日 19:21:49"
    339 = "2022年08月17日 15:06:40"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 5).Value = $updates[$row]
}
